$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("CostInvestment")
$col = $ws5.Columns.Item(1)
$col | Get-Member
